# Update crypto price/volume data per latest GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.806.31"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "2.454.90"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'517.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.86%  "
$ws.Range("D6").Value = "'132.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.60%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'0.554"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").Value = "2.461.92"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "'0.0971"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.70%  "
$ws.Range("D11").Value = "'0.156"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "'5.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "'0.335"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.33%  "
$ws.Range("D14").Value = "2.898.82"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").Value = "57.775.91"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("D16").Value = "'21.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.73%  "
$ws.Range("D17").Value = "'0.0000134"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.77%  "
$ws.Range("D18").Value = "2.464.34"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").Value = "'10.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.84%  "
$ws.Range("D20").Value = "'318.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").Value = "'4.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'5.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.00%  "
$ws.Range("D24").Value = "'64.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").Value = "'0.405"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("D26").Value = "'0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'0.159"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("D28").Value = "'7.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.72%  "
$ws.Range("D29").Value = "0.0₃0738"
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("D30").Value = "'167.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("D31").Value = "'1.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.33%  "
$ws.Range("D32").Value = "'6.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("D33").Value = "'1.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "'1.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").Value = "'17.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("D38").Value = "'3.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").Value = "'36.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").Value = "'1.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.81%  "
$ws.Range("D41").Value = "'0.787"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D42").Value = "'3.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("D43").Value = "'270.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.81%  "
$ws.Range("D44").Value = "'4.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.65%  "
$ws.Range("D45").Value = "'0.586"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("D46").Value = "'123.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("D47").Value = "'0.0902"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").Value = "'0.0484"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("D49").Value = "'0.0211"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.89%  "
$ws.Range("D50").Value = "'16.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.79%  "
$ws.Range("D51").Value = "1.718.33"
$ws.Range("E51").Value = "  -1.94%  "
